# Updated cryptos list values (price + 1h volume) per the source diff.
# Numeric-looking Price values are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.444.42'
$ws.Range("E2").Value = '  -2.12%  '
$ws.Range("D3").Value = '1.843.47'
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("D4").Value = "'" + '1.000'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'" + '260.32'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").Value = "'" + '0.3238'
$ws.Range("E8").Value = '  -8.41%  '
$ws.Range("D9").Value = "'" + '0.06744'
$ws.Range("E9").Value = '  -4.16%  '
$ws.Range("D10").Value = "'" + '18.91'
$ws.Range("E10").Value = '  -6.90%  '
$ws.Range("D11").Value = "'" + '0.7713'
$ws.Range("E11").Value = '  -5.41%  '
$ws.Range("D12").Value = "'" + '0.07687'
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("D13").Value = '1.855.24'
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("D14").Value = "'" + '89.12'
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = "'" + '5.024'
$ws.Range("D16").Value = "'" + '1.001'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = "'" + '14.12'
$ws.Range("E17").Value = '  -3.27%  '
$ws.Range("D18").Value = "'" + '1.000'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").Value = "'" + '0.000007879'
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("D20").Value = '26.476.72'
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("D21").Value = '2.079.21'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").Value = "'" + '4.532'
$ws.Range("E22").Value = '  -4.97%  '
$ws.Range("D23").Value = "'" + '9.462'
$ws.Range("E23").Value = '  -7.13%  '
$ws.Range("D24").Value = "'" + '5.909'
$ws.Range("E24").Value = '  -5.21%  '
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = "'" + '144.21'
$ws.Range("E26").Value = '  -1.70%  '
$ws.Range("D27").Value = "'" + '1.646'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").Value = "'" + '16.87'
$ws.Range("E28").Value = '  -4.05%  '
$ws.Range("D29").Value = "'" + '111.30'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("D30").Value = "'" + '4.181'
$ws.Range("E30").Value = '  -4.47%  '
$ws.Range("D31").Value = "'" + '0.08769'
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").Value = "'" + '4.124'
$ws.Range("E32").Value = '  -5.79%  '
$ws.Range("D33").Value = "'" + '0.04843'
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").Value = "'" + '1.132'
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").Value = "'" + '2.850'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = "'" + '0.6859'
$ws.Range("E36").Value = '  -7.55%  '
$ws.Range("D37").Value = "'" + '3.110'
$ws.Range("E37").Value = '  -5.83%  '
$ws.Range("D38").Value = "'" + '0.01789'
$ws.Range("E38").Value = '  -4.95%  '
$ws.Range("D39").Value = "'" + '2.219'
$ws.Range("E39").Value = '  -7.95%  '
$ws.Range("D40").Value = "'" + '0.4915'
$ws.Range("E40").Value = '  -7.27%  '
$ws.Range("D41").Value = "'" + '112.83'
$ws.Range("E41").Value = '  -3.64%  '
$ws.Range("D42").Value = "'" + '0.8983'
$ws.Range("E42").Value = '  -8.45%  '
$ws.Range("D43").Value = "'" + '6.175'
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("D44").Value = "'" + '1.0000'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = "'" + '7.754'
$ws.Range("E45").Value = '  -5.10%  '
$ws.Range("D46").Value = "'" + '0.4196'
$ws.Range("E46").Value = '  -8.81%  '
$ws.Range("D47").Value = "'" + '0.1259'
$ws.Range("E47").Value = '  -7.74%  '
$ws.Range("D48").Value = "'" + '9.083'
$ws.Range("E48").Value = '  -3.95%  '
$ws.Range("D49").Value = "'" + '0.05877'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").Value = "'" + '35.44'
$ws.Range("E50").Value = '  -3.48%  '
$ws.Range("D51").Value = "'" + '59.27'
$ws.Range("E51").Value = '  -4.11%  '